$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------------
# Two new data rows are introduced further down (the sheet grows from
# A1:N77 to A1:N79). Insert two blank rows at the right spots so rows
# currently at 38.. ("Over Old Hills", "Long Lost to Where No Pathway Goes",
# ...) shift down into their new homes (39, 41, 42, ...), while rows 36-37
# stay put.
# --------------------------------------------------------------------------
$ws.Rows("38:38").Insert()
$ws.Rows("40:40").Insert()

# Row 36 currently carries a one-off custom row style (s="7", customFormat)
# that the target workbook no longer has -- strip it, keeping its values.
$ws.Rows("36:36").ClearFormats()

# Row 36: "Kôr" entry gains a location/place excerpt.
$ws.Range("D36").Value2 = "All thy trees, Kôrtirion, were bent, since first the elves here built ancient, renowned Kôrtirion"
$ws.Range("E36").Value2 = "Kôrtirion"

# Row 37: "Unto a Long Glory..." entry gains an Erech reference.
$ws.Range("D37").Value2 = "At the stone of Erech they shall stand again"
$ws.Range("E37").Value2 = "Erech"

# Row 38: brand-new duplicate "Unto a Long Glory..." entry, now dated 1998,
# referencing the Paths of the Dead.
$ws.Range("A38").Value2 = 1998
$ws.Range("B38").Value2 = "Dol Guldur"
$ws.Range("C38").Value2 = "Unto a Long Glory..."
$ws.Range("D38").Value2 = "To the Paths of the Dead"
$ws.Range("E38").Value2 = "Paths of the Dead"

# Row 39: the former row 38 ("Over Old Hills") now carries the
# Cottage of Lost Play / Tol Eressëa annotations.
$ws.Range("D39").Value2 = "The Cottage of Lost Play"
$ws.Range("E39").Value2 = "Tol Eressëa"
$ws.Range("F39").Value2 = "The Lonely Island, The Lonely Isle, Eressëa"
$ws.Range("H39").Value2 = "island"
$ws.Range("J39").Value2 = "Off the east coast of Aman in the Bay of Eldamar"
$ws.Range("N39").Value2 = "https://lotr.fandom.com/wiki/Cottage_of_Lost_Play"

# Row 40: brand-new standalone Rhûn entry (no year/album columns).
$ws.Range("C40").Value2 = "Rhûn"
$ws.Range("D40").Value2 = "Rhûn"
$ws.Range("F40").Value2 = "The East, Eastlands"

# Match the recorded selection/scroll state of the saved file.
$ws.Application.Goto($ws.Range("A29"))
$ws.Range("G56").Select()
